# Update automàtic: dades i banners [2026-02-13 06:50]
#
# Refreshes the MeteoCat daily-summary extraction columns (DATA_EXTRACCIO,
# HUMITAT_MITJANA_DIA, PRECIPITACIO_ACUM_DIA, PRESSIO_ATMOSFERICA,
# RADIACIO_GLOBAL, RATXA_VENT_MAX, TEMPERATURA_*_DIA, GRUIX_NEU_MAX) with the
# 06:50 scrape values for every station row (2-46).
#
# Percentage-shaped text (e.g. "68%") would otherwise be auto-coerced by
# Excel into a numeric percentage (changing both the stored value and the
# cell's number-format style) when assigned straight to .Value/.Formula.
# To keep these columns as plain text - matching the sheet's existing
# t="inlineStr" string cells and untouched style index - we route them
# through a TEXT() formula and then Copy/PasteSpecial(xlPasteValues) to
# freeze the formula result back down to a literal string in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

function Set-TextValue($range, [string]$value) {
    if ($value.TrimEnd().EndsWith('%')) {
        # Force text so Excel doesn't reinterpret "NN%" as a numeric percentage.
        $escaped = $value.Replace('"', '""')
        $range.Formula = '=TEXT("' + $escaped + '","@")'
        $range.Copy()
        $range.PasteSpecial($xlPasteValues)
    } else {
        $range.Value = $value
    }
}

# Row 2
Set-TextValue $ws.Range('E2') '2026-02-13 06:48:48'
Set-TextValue $ws.Range('O2') '-2.2 °C'

# Row 3
Set-TextValue $ws.Range('E3') '2026-02-13 06:48:51'
Set-TextValue $ws.Range('H3') '68%'
Set-TextValue $ws.Range('I3') '0.1 mm'

# Row 4
Set-TextValue $ws.Range('E4') '2026-02-13 06:48:53'
Set-TextValue $ws.Range('H4') '58%'
Set-TextValue $ws.Range('J4') '1002.5 hPa'
Set-TextValue $ws.Range('N4') '8.0 °C 6:16 TU'
Set-TextValue $ws.Range('O4') '9.6 °C'

# Row 5
Set-TextValue $ws.Range('E5') '2026-02-13 06:48:56'
Set-TextValue $ws.Range('H5') '67%'
Set-TextValue $ws.Range('L5') '37.8 km/h - 89º 6:13 TU'

# Row 6
Set-TextValue $ws.Range('E6') '2026-02-13 06:48:59'
Set-TextValue $ws.Range('J6') '1002.7 hPa'

# Row 7
Set-TextValue $ws.Range('E7') '2026-02-13 06:49:01'
Set-TextValue $ws.Range('J7') '1003.1 hPa'
Set-TextValue $ws.Range('O7') '14.4 °C'

# Row 8
Set-TextValue $ws.Range('E8') '2026-02-13 06:49:04'
Set-TextValue $ws.Range('J8') '1003.0 hPa'

# Row 9
Set-TextValue $ws.Range('E9') '2026-02-13 06:49:07'
Set-TextValue $ws.Range('H9') '67%'
Set-TextValue $ws.Range('O9') '8.3 °C'

# Row 10
Set-TextValue $ws.Range('E10') '2026-02-13 06:49:10'
Set-TextValue $ws.Range('H10') '79%'
Set-TextValue $ws.Range('O10') '7.6 °C'

# Row 11
Set-TextValue $ws.Range('E11') '2026-02-13 06:49:12'

# Row 12
Set-TextValue $ws.Range('E12') '2026-02-13 06:49:15'
Set-TextValue $ws.Range('O12') '7.9 °C'

# Row 13
Set-TextValue $ws.Range('E13') '2026-02-13 06:49:17'
Set-TextValue $ws.Range('J13') '1006.6 hPa'

# Row 14
Set-TextValue $ws.Range('E14') '2026-02-13 06:49:20'
Set-TextValue $ws.Range('H14') '63%'
Set-TextValue $ws.Range('O14') '11.3 °C'

# Row 15
Set-TextValue $ws.Range('E15') '2026-02-13 06:49:23'
Set-TextValue $ws.Range('H15') '69%'
Set-TextValue $ws.Range('O15') '8.5 °C'

# Row 16
Set-TextValue $ws.Range('E16') '2026-02-13 06:49:26'
Set-TextValue $ws.Range('H16') '63%'
Set-TextValue $ws.Range('N16') '-5.3 °C 6:24 TU'
Set-TextValue $ws.Range('O16') '-3.1 °C'

# Row 17
Set-TextValue $ws.Range('E17') '2026-02-13 06:49:28'
Set-TextValue $ws.Range('H17') '69%'

# Row 18
Set-TextValue $ws.Range('E18') '2026-02-13 06:49:31'
Set-TextValue $ws.Range('J18') '1002.8 hPa'
Set-TextValue $ws.Range('O18') '7.4 °C'

# Row 19
Set-TextValue $ws.Range('E19') '2026-02-13 06:49:34'
Set-TextValue $ws.Range('H19') '71%'

# Row 20
Set-TextValue $ws.Range('E20') '2026-02-13 06:49:36'
Set-TextValue $ws.Range('H20') '80%'
Set-TextValue $ws.Range('I20') '0.4 mm'

# Row 21
Set-TextValue $ws.Range('E21') '2026-02-13 06:49:39'
Set-TextValue $ws.Range('J21') '1005.2 hPa'

# Row 22
Set-TextValue $ws.Range('E22') '2026-02-13 06:49:42'
Set-TextValue $ws.Range('G22') '116 cm'
Set-TextValue $ws.Range('H22') '84%'
Set-TextValue $ws.Range('I22') '1.0 mm'

# Row 23
Set-TextValue $ws.Range('E23') '2026-02-13 06:49:45'
Set-TextValue $ws.Range('H23') '68%'

# Row 24
Set-TextValue $ws.Range('E24') '2026-02-13 06:49:48'
Set-TextValue $ws.Range('H24') '86%'
Set-TextValue $ws.Range('I24') '1.3 mm'
Set-TextValue $ws.Range('J24') '1004.0 hPa'

# Row 25
Set-TextValue $ws.Range('E25') '2026-02-13 06:49:51'
Set-TextValue $ws.Range('H25') '60%'
Set-TextValue $ws.Range('N25') '-4.5 °C 6:29 TU'
Set-TextValue $ws.Range('O25') '-2.8 °C'

# Row 26
Set-TextValue $ws.Range('E26') '2026-02-13 06:49:53'
Set-TextValue $ws.Range('H26') '54%'
Set-TextValue $ws.Range('J26') '1003.2 hPa'
Set-TextValue $ws.Range('K26') '-0.1 MJ/m2'
Set-TextValue $ws.Range('O26') '2.5 °C'

# Row 27
Set-TextValue $ws.Range('E27') '2026-02-13 06:49:56'
Set-TextValue $ws.Range('H27') '63%'
Set-TextValue $ws.Range('I27') '0.7 mm'
Set-TextValue $ws.Range('N27') '-3.3 °C 6:23 TU'
Set-TextValue $ws.Range('O27') '-1.7 °C'

# Row 28
Set-TextValue $ws.Range('E28') '2026-02-13 06:49:59'
Set-TextValue $ws.Range('H28') '68%'
Set-TextValue $ws.Range('J28') '1003.2 hPa'

# Row 29
Set-TextValue $ws.Range('E29') '2026-02-13 06:50:02'
Set-TextValue $ws.Range('O29') '10.4 °C'

# Row 30
Set-TextValue $ws.Range('E30') '2026-02-13 06:50:04'
Set-TextValue $ws.Range('H30') '75%'
Set-TextValue $ws.Range('J30') '1002.9 hPa'

# Row 31
Set-TextValue $ws.Range('E31') '2026-02-13 06:50:07'
Set-TextValue $ws.Range('H31') '56%'
Set-TextValue $ws.Range('J31') '1002.0 hPa'
Set-TextValue $ws.Range('N31') '10.5 °C 6:10 TU'
Set-TextValue $ws.Range('O31') '11.5 °C'

# Row 32
Set-TextValue $ws.Range('E32') '2026-02-13 06:50:10'
Set-TextValue $ws.Range('K32') '-0.1 MJ/m2'

# Row 33
Set-TextValue $ws.Range('E33') '2026-02-13 06:50:13'
Set-TextValue $ws.Range('J33') '1005.3 hPa'

# Row 34
Set-TextValue $ws.Range('E34') '2026-02-13 06:50:15'
Set-TextValue $ws.Range('H34') '56%'
Set-TextValue $ws.Range('O34') '-0.5 °C'

# Row 35
Set-TextValue $ws.Range('E35') '2026-02-13 06:50:18'
Set-TextValue $ws.Range('H35') '57%'
Set-TextValue $ws.Range('J35') '1004.5 hPa'
Set-TextValue $ws.Range('O35') '6.4 °C'

# Row 36
Set-TextValue $ws.Range('E36') '2026-02-13 06:50:21'
Set-TextValue $ws.Range('H36') '62%'
Set-TextValue $ws.Range('J36') '1002.7 hPa'
Set-TextValue $ws.Range('O36') '11.0 °C'

# Row 37
Set-TextValue $ws.Range('E37') '2026-02-13 06:50:24'
Set-TextValue $ws.Range('H37') '65%'
Set-TextValue $ws.Range('J37') '1004.8 hPa'
Set-TextValue $ws.Range('O37') '3.7 °C'

# Row 38
Set-TextValue $ws.Range('E38') '2026-02-13 06:50:27'
Set-TextValue $ws.Range('H38') '55%'
Set-TextValue $ws.Range('O38') '10.2 °C'

# Row 39
Set-TextValue $ws.Range('E39') '2026-02-13 06:50:29'
Set-TextValue $ws.Range('H39') '53%'
Set-TextValue $ws.Range('I39') '0.1 mm'
Set-TextValue $ws.Range('N39') '-5.4 °C 6:29 TU'
Set-TextValue $ws.Range('O39') '-2.8 °C'

# Row 40
Set-TextValue $ws.Range('E40') '2026-02-13 06:50:32'
Set-TextValue $ws.Range('J40') '1006.3 hPa'

# Row 41
Set-TextValue $ws.Range('E41') '2026-02-13 06:50:35'
Set-TextValue $ws.Range('H41') '52%'
Set-TextValue $ws.Range('J41') '1003.4 hPa'
Set-TextValue $ws.Range('N41') '10.0 °C 6:00 TU'
Set-TextValue $ws.Range('O41') '12.7 °C'

# Row 42
Set-TextValue $ws.Range('E42') '2026-02-13 06:50:38'
Set-TextValue $ws.Range('O42') '10.8 °C'

# Row 43
Set-TextValue $ws.Range('E43') '2026-02-13 06:50:40'
Set-TextValue $ws.Range('H43') '67%'
Set-TextValue $ws.Range('O43') '7.1 °C'

# Row 44
Set-TextValue $ws.Range('E44') '2026-02-13 06:50:43'
Set-TextValue $ws.Range('H44') '82%'

# Row 45
Set-TextValue $ws.Range('E45') '2026-02-13 06:50:46'
Set-TextValue $ws.Range('H45') '67%'
Set-TextValue $ws.Range('J45') '1003.0 hPa'
Set-TextValue $ws.Range('K45') '-0.1 MJ/m2'
Set-TextValue $ws.Range('L45') '28.4 km/h - 147º 6:19 TU'
Set-TextValue $ws.Range('M45') '7.9 °C 6:21 TU'
Set-TextValue $ws.Range('O45') '3.0 °C'

# Row 46
Set-TextValue $ws.Range('E46') '2026-02-13 06:50:49'
Set-TextValue $ws.Range('I46') '0.2 mm'
Set-TextValue $ws.Range('J46') '1004.3 hPa'
Set-TextValue $ws.Range('L46') '10.1 km/h - 165º 6:12 TU'

$excel.CutCopyMode = $false
